$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_WVR = $wb.Worksheets.Item("WVR")

# Row 2 (ALC)
$ws_ALC.Range("H2").Value = 208
$ws_ALC.Range("I2").Value = 212.5
$ws_ALC.Range("J2").Value = 190
$ws_ALC.Range("K2").Value = 212.5
$ws_ALC.Range("L2").Value = 190
$ws_ALC.Range("M2").Value = -99.5
$ws_ALC.Range("N2").Value = -416

# Row 33 (ALC)
$ws_ALC.Range("H33").Value = 30303642
$ws_ALC.Range("I33").Value = 597.8261
$ws_ALC.Range("K33").Value = 597.8261
$ws_ALC.Range("M33").Value = -368.8261

# Row 61 (ALC)
$ws_ALC.Range("H61").Value = 215.8
$ws_ALC.Range("I61").Value = 215.8
$ws_ALC.Range("K61").Value = 647.4000000000001
$ws_ALC.Range("M61").Value = -475.4000000000001

# Row 76 (ALC)
$ws_ALC.Range("H76").Value = 3041.6667
$ws_ALC.Range("I76").Value = 3028.5715
$ws_ALC.Range("J76").Value = 3133.3333
$ws_ALC.Range("K76").Value = 3028.5715
$ws_ALC.Range("L76").Value = 3133.3333
$ws_ALC.Range("M76").Value = -2713.5715
$ws_ALC.Range("N76").Value = -3763.3333

# Row 79 (ALC)
$ws_ALC.Range("H79").Value = 3041.6667
$ws_ALC.Range("I79").Value = 3028.5715
$ws_ALC.Range("J79").Value = 3133.3333
$ws_ALC.Range("K79").Value = 3028.5715
$ws_ALC.Range("L79").Value = 3133.3333
$ws_ALC.Range("M79").Value = -1936.5715
$ws_ALC.Range("N79").Value = -5317.3333

# Row 115 (ALC)
$ws_ALC.Range("H115").Value = 724.5
$ws_ALC.Range("I115").Value = 416.1111
$ws_ALC.Range("J115").Value = 3500
$ws_ALC.Range("K115").Value = 1248.3333
$ws_ALC.Range("L115").Value = 10500
$ws_ALC.Range("M115").Value = 318.6667
$ws_ALC.Range("N115").Value = -13634

# Row 116 (ALC)
$ws_ALC.Range("H116").Value = 5914.643
$ws_ALC.Range("I116").Value = 2499.8
$ws_ALC.Range("J116").Value = 14451.75
$ws_ALC.Range("K116").Value = 2499.8
$ws_ALC.Range("L116").Value = 14451.75
$ws_ALC.Range("M116").Value = 942.1999999999998
$ws_ALC.Range("N116").Value = -21335.75

# Row 118 (ALC)
$ws_ALC.Range("H118").Value = 603.5294
$ws_ALC.Range("I118").Value = 251.11111
$ws_ALC.Range("K118").Value = 753.3333299999999
$ws_ALC.Range("M118").Value = 903.6666700000001

# Row 123 (ALC)
$ws_ALC.Range("H123").Value = 26666.666
$ws_ALC.Range("J123").Value = 26666.666
$ws_ALC.Range("L123").Value = 26666.666
$ws_ALC.Range("N123").Value = -36466.666

# Row 124 (ALC)
$ws_ALC.Range("H124").Value = 25000
$ws_ALC.Range("J124").Value = 25000
$ws_ALC.Range("L124").Value = 25000
$ws_ALC.Range("N124").Value = -34820

# Row 126 (ALC)
$ws_ALC.Range("H126").Value = 24444.445
$ws_ALC.Range("J126").Value = 24444.445
$ws_ALC.Range("L126").Value = 24444.445
$ws_ALC.Range("N126").Value = -34324.445

# Row 130 (ALC)
$ws_ALC.Range("H130").Value = 30000
$ws_ALC.Range("J130").Value = 30000
$ws_ALC.Range("L130").Value = 30000
$ws_ALC.Range("N130").Value = -40040

# Row 137 (ALC)
$ws_ALC.Range("H137").Value = 2317476.2
$ws_ALC.Range("I137").Value = 2419.6316
$ws_ALC.Range("J137").Value = 11114691
$ws_ALC.Range("K137").Value = 7258.8948
$ws_ALC.Range("L137").Value = 33344073
$ws_ALC.Range("M137").Value = -4708.8948
$ws_ALC.Range("N137").Value = -33349173

# Row 55 (BSM)
$ws_BSM.Range("H55").Value = 28333
$ws_BSM.Range("J55").Value = 28333
$ws_BSM.Range("L55").Value = 28333
$ws_BSM.Range("N55").Value = -28879

# Row 31 (CRP)
$ws_CRP.Range("H31").Value = 35716188
$ws_CRP.Range("I31").Value = 100000860
$ws_CRP.Range("J31").Value = 2479.1667
$ws_CRP.Range("K31").Value = 100000860
$ws_CRP.Range("L31").Value = 2479.1667
$ws_CRP.Range("M31").Value = -100000565
$ws_CRP.Range("N31").Value = -3069.1667

# Row 34 (CRP)
$ws_CRP.Range("H34").Value = 35716188
$ws_CRP.Range("I34").Value = 100000860
$ws_CRP.Range("J34").Value = 2479.1667
$ws_CRP.Range("K34").Value = 100000860
$ws_CRP.Range("L34").Value = 2479.1667
$ws_CRP.Range("M34").Value = -100000658
$ws_CRP.Range("N34").Value = -2883.1667

# Row 99 (CRP)
$ws_CRP.Range("H99").Value = 3503.2666
$ws_CRP.Range("I99").Value = 3629.4119
$ws_CRP.Range("J99").Value = 3338.3076
$ws_CRP.Range("K99").Value = 3629.4119
$ws_CRP.Range("L99").Value = 3338.3076
$ws_CRP.Range("M99").Value = -2131.4119
$ws_CRP.Range("N99").Value = -6334.3076

# Row 105 (CRP)
$ws_CRP.Range("H105").Value = 1544.56
$ws_CRP.Range("I105").Value = 1530.75
$ws_CRP.Range("J105").Value = 1599.8
$ws_CRP.Range("K105").Value = 1530.75
$ws_CRP.Range("L105").Value = 1599.8
$ws_CRP.Range("M105").Value = 216.25
$ws_CRP.Range("N105").Value = -5093.8

# Row 126 (CRP)
$ws_CRP.Range("H126").Value = 3503.2666
$ws_CRP.Range("I126").Value = 3629.4119
$ws_CRP.Range("J126").Value = 3338.3076
$ws_CRP.Range("K126").Value = 10888.2357
$ws_CRP.Range("L126").Value = 10014.9228
$ws_CRP.Range("M126").Value = -8418.235700000001
$ws_CRP.Range("N126").Value = -14954.9228

# Row 132 (CRP)
$ws_CRP.Range("H132").Value = 2197.818
$ws_CRP.Range("I132").Value = 1954.9429
$ws_CRP.Range("J132").Value = 3142.3333
$ws_CRP.Range("K132").Value = 5864.8287
$ws_CRP.Range("L132").Value = 9426.999899999999
$ws_CRP.Range("M132").Value = -3334.8287
$ws_CRP.Range("N132").Value = -14486.9999

# Row 134 (CRP)
$ws_CRP.Range("H134").Value = 1698.0416
$ws_CRP.Range("I134").Value = 1705.65
$ws_CRP.Range("K134").Value = 5116.950000000001
$ws_CRP.Range("M134").Value = -2581.950000000001

# Row 64 (CUL)
$ws_CUL.Range("H64").Value = 4081.5
$ws_CUL.Range("I64").Value = 3168
$ws_CUL.Range("J64").Value = 5256
$ws_CUL.Range("K64").Value = 9504
$ws_CUL.Range("L64").Value = 15768
$ws_CUL.Range("M64").Value = -9234
$ws_CUL.Range("N64").Value = -16308

# Row 67 (CUL)
$ws_CUL.Range("H67").Value = 4081.5
$ws_CUL.Range("I67").Value = 3168
$ws_CUL.Range("J67").Value = 5256
$ws_CUL.Range("K67").Value = 9504
$ws_CUL.Range("L67").Value = 15768
$ws_CUL.Range("M67").Value = -8568
$ws_CUL.Range("N67").Value = -17640

# Row 131 (CUL)
$ws_CUL.Range("H131").Value = 890.8099999999999
$ws_CUL.Range("I131").Value = 868.6667
$ws_CUL.Range("J131").Value = 891.4949
$ws_CUL.Range("K131").Value = 2606.0001
$ws_CUL.Range("L131").Value = 2674.4847
$ws_CUL.Range("M131").Value = 2433.9999
$ws_CUL.Range("N131").Value = -12754.4847

# Row 70 (GSM)
$ws_GSM.Range("H70").Value = 5331.1875
$ws_GSM.Range("I70").Value = 5113.25
$ws_GSM.Range("J70").Value = 5985
$ws_GSM.Range("K70").Value = 5113.25
$ws_GSM.Range("L70").Value = 5985
$ws_GSM.Range("M70").Value = -4843.25
$ws_GSM.Range("N70").Value = -6525

# Row 73 (GSM)
$ws_GSM.Range("H73").Value = 5331.1875
$ws_GSM.Range("I73").Value = 5113.25
$ws_GSM.Range("J73").Value = 5985
$ws_GSM.Range("K73").Value = 5113.25
$ws_GSM.Range("L73").Value = 5985
$ws_GSM.Range("M73").Value = -4177.25
$ws_GSM.Range("N73").Value = -7857

# Row 107 (WVR)
$ws_WVR.Range("H107").Value = 11914.889
$ws_WVR.Range("I107").Value = 17623.5
$ws_WVR.Range("J107").Value = 497.66666
$ws_WVR.Range("K107").Value = 52870.5
$ws_WVR.Range("L107").Value = 1492.99998
$ws_WVR.Range("M107").Value = -50950.5
$ws_WVR.Range("N107").Value = -5332.999980000001

# Row 126 (WVR)
$ws_WVR.Range("H126").Value = 62501956
$ws_WVR.Range("I126").Value = 71430100
$ws_WVR.Range("J126").Value = 5002.5
$ws_WVR.Range("K126").Value = 214290300
$ws_WVR.Range("L126").Value = 15007.5
$ws_WVR.Range("M126").Value = -214287830
$ws_WVR.Range("N126").Value = -19947.5
